$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OPEX")

$ws.Range("B3").Value = 19283.042118970367
$ws.Range("B4").Value = 43709.823692218881
$ws.Range("B5").Value = 44248.4220525376
$ws.Range("B6").Value = 10540.776519765759
$ws.Range("B7").Value = 11938.918689039821
$ws.Range("B8").Value = 50983.71612182937
$ws.Range("B9").Value = 47726.863400003029
$ws.Range("B10").Value = 48220.289169209609
$ws.Range("B11").Value = 37682.867713520049
$ws.Range("B12").Value = 44827.8721002112
$ws.Range("B13").Value = 12139.814930199653
$ws.Range("B14").Value = 47419.574760442287
$ws.Range("B15").Value = 50088.476920599649
